$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.31"
$ws.Range("E2").Value = "'0.64%"
$ws.Range("D3").Value = "'31.49"
$ws.Range("E3").Value = "'0.92%"
$ws.Range("D4").Value = "'5.074"
$ws.Range("E4").Value = "'-1.29%"
$ws.Range("D5").Value = "'0.07845"
$ws.Range("E5").Value = "'-2.49%"
$ws.Range("D6").Value = "'2.346"
$ws.Range("E6").Value = "'-11.22%"
$ws.Range("D7").Value = "'7.816"
$ws.Range("E7").Value = "'-0.47%"
$ws.Range("D8").Value = "'3.833"
$ws.Range("E8").Value = "'-0.01%"
$ws.Range("D9").Value = "'0.9182"
$ws.Range("E9").Value = "'1.06%"
$ws.Range("E10").Value = "'1.83%"
$ws.Range("D11").Value = "'0.07553"
$ws.Range("D12").Value = "'0.09163"
$ws.Range("E12").Value = "'13.82%"
$ws.Range("D13").Value = "'0.02991"
$ws.Range("E13").Value = "'-1.34%"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.001507"
$ws.Range("E15").Value = "'0.37%"
$ws.Range("D16").Value = "'0.005901"
$ws.Range("E16").Value = "'-0.96%"
$ws.Range("E17").Value = "'-0.96%"
$ws.Range("E19").Value = "'-0.52%"
$ws.Range("D20").Value = "'0.1337"
$ws.Range("E20").Value = "'0.55%"
$ws.Range("D21").Value = "'4.007"
$ws.Range("E21").Value = "'-12.84%"
$ws.Range("E22").Value = "'11.70%"
$ws.Range("D23").Value = "'0.04618"
$ws.Range("E23").Value = "'0.23%"
$ws.Range("D24").Value = "'0.001251"
$ws.Range("E24").Value = "'-0.84%"
$ws.Range("D25").Value = "'0.004472"
$ws.Range("E25").Value = "'0.65%"
$ws.Range("D26").Value = "'0.0001249"
$ws.Range("E26").Value = "'5.74%"
$ws.Range("E27").Value = "'-1.53%"
$ws.Range("D39").Value = "'0.01764"
$ws.Range("E39").Value = "'-2.28%"
$ws.Range("D40").Value = "'0.04760"
$ws.Range("E40").Value = "'4.89%"
$ws.Range("D41").Value = "'0.007160"
$ws.Range("E41").Value = "'0.85%"
$ws.Range("D42").Value = "'0.1361"
$ws.Range("E42").Value = "'1.43%"
$ws.Range("D43").Value = "'0.002188"
$ws.Range("E43").Value = "'-2.42%"
$ws.Range("D44").Value = "'0.01030"
$ws.Range("E44").Value = "'-1.27%"
$ws.Range("D45").Value = "'0.00006318"
$ws.Range("E45").Value = "'0.06%"
$ws.Range("E46").Value = "'-0.12%"
$ws.Range("E47").Value = "'28.64%"
$ws.Range("D48").Value = "'0.7431"
$ws.Range("E48").Value = "'-9.44%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.12%"
